$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4939.411
$ws.Range("I15").Value = 4939.411
$ws.Range("K15").Value = 14818.233
$ws.Range("M15").Value = -14649.233
$ws.Range("H17").Value = 565.8687
$ws.Range("I17").Value = 110
$ws.Range("J17").Value = 585.0632000000001
$ws.Range("K17").Value = 330
$ws.Range("L17").Value = 1755.1896
$ws.Range("M17").Value = -162
$ws.Range("N17").Value = -2091.1896
$ws.Range("H132").Value = 116212.23
$ws.Range("I132").Value = 123237.375
$ws.Range("J132").Value = 999.8
$ws.Range("K132").Value = 369712.125
$ws.Range("L132").Value = 2999.4
$ws.Range("M132").Value = -367182.125
$ws.Range("N132").Value = -8059.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5525.2104
$ws.Range("I32").Value = 5859.6
$ws.Range("J32").Value = 4271.25
$ws.Range("K32").Value = 5859.6
$ws.Range("L32").Value = 4271.25
$ws.Range("M32").Value = -5572.6
$ws.Range("N32").Value = -4845.25
$ws.Range("H45").Value = 830.3333
$ws.Range("I45").Value = 735.8461
$ws.Range("J45").Value = 1076
$ws.Range("K45").Value = 735.8461
$ws.Range("L45").Value = 1076
$ws.Range("M45").Value = -358.8461
$ws.Range("N45").Value = -1830
$ws.Range("H61").Value = 1447.5625
$ws.Range("I61").Value = 1429.5
$ws.Range("J61").Value = 1501.75
$ws.Range("K61").Value = 1429.5
$ws.Range("L61").Value = 1501.75
$ws.Range("M61").Value = -1217.5
$ws.Range("N61").Value = -1925.75
$ws.Range("H123").Value = 25163.334
$ws.Range("J123").Value = 25163.334
$ws.Range("L123").Value = 25163.334
$ws.Range("N123").Value = -34963.334
$ws.Range("H136").Value = 1447.5625
$ws.Range("I136").Value = 1429.5
$ws.Range("J136").Value = 1501.75
$ws.Range("K136").Value = 4288.5
$ws.Range("L136").Value = 4505.25
$ws.Range("M136").Value = -1738.5
$ws.Range("N136").Value = -9605.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1903.909
$ws.Range("I134").Value = 2098.6785
$ws.Range("J134").Value = 813.2
$ws.Range("K134").Value = 6296.0355
$ws.Range("L134").Value = 2439.6
$ws.Range("M134").Value = -3761.0355
$ws.Range("N134").Value = -7509.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1027.2727
$ws.Range("I5").Value = 1600
$ws.Range("K5").Value = 4800
$ws.Range("M5").Value = -4688
$ws.Range("H50").Value = 202.9
$ws.Range("I50").Value = 214.33333
$ws.Range("J50").Value = 100
$ws.Range("K50").Value = 642.99999
$ws.Range("L50").Value = 300
$ws.Range("M50").Value = -161.99999
$ws.Range("N50").Value = -1262
$ws.Range("H53").Value = 202.9
$ws.Range("I53").Value = 214.33333
$ws.Range("J53").Value = 100
$ws.Range("K53").Value = 642.99999
$ws.Range("L53").Value = 300
$ws.Range("M53").Value = -161.99999
$ws.Range("N53").Value = -1262
$ws.Range("H131").Value = 4846.2
$ws.Range("I131").Value = 6620.125
$ws.Range("J131").Value = 4011.4119
$ws.Range("K131").Value = 19860.375
$ws.Range("L131").Value = 12034.2357
$ws.Range("M131").Value = -14820.375
$ws.Range("N131").Value = -22114.2357
$ws.Range("H133").Value = 9088.888999999999
$ws.Range("I133").Value = 5000
$ws.Range("K133").Value = 15000
$ws.Range("M133").Value = -9940
$ws.Range("H135").Value = 1027.2727
$ws.Range("I135").Value = 1600
$ws.Range("K135").Value = 14400
$ws.Range("M135").Value = -11865

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5353.077
$ws.Range("J70").Value = 5538
$ws.Range("L70").Value = 5538
$ws.Range("N70").Value = -6078
$ws.Range("H73").Value = 5353.077
$ws.Range("J73").Value = 5538
$ws.Range("L73").Value = 5538
$ws.Range("N73").Value = -7410
$ws.Range("H80").Value = 9500
$ws.Range("I80").Value = 2628.5715
$ws.Range("J80").Value = 14844.444
$ws.Range("K80").Value = 2628.5715
$ws.Range("L80").Value = 14844.444
$ws.Range("M80").Value = -1630.5715
$ws.Range("N80").Value = -16840.444
$ws.Range("H83").Value = 9500
$ws.Range("I83").Value = 2628.5715
$ws.Range("J83").Value = 14844.444
$ws.Range("K83").Value = 13142.8575
$ws.Range("L83").Value = 74222.22
$ws.Range("M83").Value = -8150.8575
$ws.Range("N83").Value = -84206.22

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 608.5
$ws.Range("I22").Value = 518
$ws.Range("J22").Value = 699
$ws.Range("K22").Value = 518
$ws.Range("L22").Value = 699
$ws.Range("M22").Value = -223
$ws.Range("N22").Value = -1289
$ws.Range("H27").Value = 608.5
$ws.Range("I27").Value = 518
$ws.Range("J27").Value = 699
$ws.Range("K27").Value = 518
$ws.Range("L27").Value = 699
$ws.Range("M27").Value = -411
$ws.Range("N27").Value = -913
$ws.Range("H40").Value = 1213.3846
$ws.Range("I40").Value = 1260.8182
$ws.Range("J40").Value = 952.5
$ws.Range("K40").Value = 1260.8182
$ws.Range("L40").Value = 952.5
$ws.Range("M40").Value = -1124.8182
$ws.Range("N40").Value = -1224.5
$ws.Range("H46").Value = 2280.2
$ws.Range("I46").Value = 1450
$ws.Range("J46").Value = 2833.6667
$ws.Range("K46").Value = 1450
$ws.Range("L46").Value = 2833.6667
$ws.Range("M46").Value = -1262
$ws.Range("N46").Value = -3209.6667
$ws.Range("H61").Value = 1143.2069
$ws.Range("I61").Value = 1041.5238
$ws.Range("J61").Value = 1410.125
$ws.Range("K61").Value = 1041.5238
$ws.Range("L61").Value = 1410.125
$ws.Range("M61").Value = -839.5237999999999
$ws.Range("N61").Value = -1814.125
$ws.Range("H68").Value = 1902.52
$ws.Range("I68").Value = 1816.875
$ws.Range("J68").Value = 2054.7778
$ws.Range("K68").Value = 1816.875
$ws.Range("L68").Value = 2054.7778
$ws.Range("M68").Value = -1067.875
$ws.Range("N68").Value = -3552.7778
$ws.Range("H71").Value = 1902.52
$ws.Range("I71").Value = 1816.875
$ws.Range("J71").Value = 2054.7778
$ws.Range("K71").Value = 9084.375
$ws.Range("L71").Value = 10273.889
$ws.Range("M71").Value = -5340.375
$ws.Range("N71").Value = -17761.889
$ws.Range("H82").Value = 1578.091
$ws.Range("I82").Value = 1211.2858
$ws.Range("K82").Value = 1211.2858
$ws.Range("M82").Value = -850.2858000000001
$ws.Range("H85").Value = 1578.091
$ws.Range("I85").Value = 1211.2858
$ws.Range("K85").Value = 1211.2858
$ws.Range("M85").Value = 36.71419999999989
$ws.Range("H113").Value = 1143.2069
$ws.Range("I113").Value = 1041.5238
$ws.Range("J113").Value = 1410.125
$ws.Range("K113").Value = 1041.5238
$ws.Range("L113").Value = 1410.125
$ws.Range("M113").Value = 1128.4762
$ws.Range("N113").Value = -5750.125
$ws.Range("H122").Value = 1884.1538
$ws.Range("I122").Value = 998
$ws.Range("J122").Value = 2150
$ws.Range("K122").Value = 2994
$ws.Range("L122").Value = 6450
$ws.Range("M122").Value = -544
$ws.Range("N122").Value = -11350
$ws.Range("H127").Value = 34950
$ws.Range("J127").Value = 34950
$ws.Range("L127").Value = 34950
$ws.Range("N127").Value = -44870
$ws.Range("H136").Value = 2848.5881
$ws.Range("I136").Value = 1066.1428
$ws.Range("K136").Value = 3198.4284
$ws.Range("M136").Value = -648.4284000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 29289.166
$ws.Range("I62").Value = 34828.332
$ws.Range("J62").Value = 23750
$ws.Range("K62").Value = 34828.332
$ws.Range("L62").Value = 23750
$ws.Range("M62").Value = -34204.332
$ws.Range("N62").Value = -24998
$ws.Range("H65").Value = 29289.166
$ws.Range("I65").Value = 34828.332
$ws.Range("J65").Value = 23750
$ws.Range("K65").Value = 174141.66
$ws.Range("L65").Value = 118750
$ws.Range("M65").Value = -171021.66
$ws.Range("N65").Value = -124990
